# Add data for 2022-05-26
# - Advance the "through" date in the sheet name and header text from May 17 to May 18
# - Update historical "May" counts (columns B, G, L, Q, V, AA, AF, AK, ...) for several
#   neighborhoods as new/late-reported incidents are reflected in the data
# - Correct one misclassified 2022 record (New City -> Chinatown)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab
$ws.Name = "Through 2022-05-18"

# 2) Update the column header text in B1 (shared string) to match the new "through" date
$ws.Range("B1").Value = "May 2022 (through May 18)"

# 3) Cell value updates (existing values incremented)
$ws.Range("AA4").Value = 2    # Humboldt Park, May 2017: 1 -> 2
$ws.Range("G8").Value = 3     # South Shore, May 2021: 2 -> 3
$ws.Range("AA20").Value = 2   # Woodlawn, May 2017: 1 -> 2
$ws.Range("G23").Value = 3    # Grand Crossing, May 2021: 2 -> 3
$ws.Range("AF29").Value = 2   # West Pullman, May 2016: 1 -> 2
$ws.Range("B57").Value = 2    # Chinatown, May 2022: 1 -> 2

# 4) New cell values (previously empty)
$ws.Range("G17").Value = 1    # South Chicago, May 2021
$ws.Range("L23").Value = 1    # Grand Crossing, May 2020
$ws.Range("AF25").Value = 1   # Auburn Gresham, May 2016
$ws.Range("Q28").Value = 1    # West Town, May 2019
$ws.Range("AF28").Value = 1   # West Town, May 2016
$ws.Range("Q34").Value = 1    # South Deering, May 2019
$ws.Range("L50").Value = 1    # Archer Heights, May 2020
$ws.Range("G67").Value = 1    # Hegewisch, May 2021

# 5) Remove an incorrect value (New City, May 2022 reclassified to Chinatown above)
$ws.Range("B39").ClearContents()
